$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Data Type column: CHAR -> VARCHAR for the string/text fields
# (drug_name, drug_generic name, state_name, dr_type)
$ws.Range("C4").Value2  = "VARCHAR"
$ws.Range("C5").Value2  = "VARCHAR"
$ws.Range("C10").Value2 = "VARCHAR"
$ws.Range("C15").Value2 = "VARCHAR"

# Remove the blank spacer row (old row 17) so the "primary key" / "foreign key"
# rows move up from 18/19 to 17/18
$ws.Rows(17).Delete()

# Update the view: zoom to 126% and select the used range A1:F18
$excel.ActiveWindow.Zoom = 126
$ws.Range("A1:F18").Select() | Out-Null
